# BIS-1630: Fixed xls export tests
# Insert a new "Internal" column before column B, shifting existing
# columns B..N to C..O, and add a new "Internal Assignment" column at
# the end (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the new column before column B. This shifts B:N -> C:O
#    and carries over values/styles/types unchanged.
$ws.Columns("B:B").Insert()

# 2) Header labels for the new "Internal" column.
$ws.Range("B2").Value = "Internal"
$ws.Range("B4").Value = "Internal"

# 3) Property-code text values lost their leading "$" prefix.
$ws.Range("A5").Value = "NAME"
$ws.Range("A6").Value = "DOCUMENT"
$ws.Range("A7").Value = "ANNOTATIONS_STATE"

# 4) New "Internal" flag values for the data rows. Copy from existing
#    text-typed TRUE/FALSE cells so the new cells keep text type
#    (matching the exported shared-string cells) instead of being
#    auto-coerced to booleans.
$ws.Range("C5").Copy($ws.Range("B3"))
$ws.Range("D3").Copy($ws.Range("B5"))
$ws.Range("D3").Copy($ws.Range("B6"))
$ws.Range("D3").Copy($ws.Range("B7"))

# 5) New trailing "Internal Assignment" column.
$ws.Range("P4").Value = "Internal Assignment"
$ws.Range("C5").Copy($ws.Range("P5"))
$ws.Range("C5").Copy($ws.Range("P6"))
$ws.Range("C5").Copy($ws.Range("P7"))

# 6) Sheet view / cursor bookkeeping to match the saved selection state.
$ws.Range("A1").Select()
$ws.Range("P4:P7").Select()
